$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7; this shifts the existing rows 7..73 down
# to 8..74 (their data is carried along automatically), matching the
# dimension growing from A1:R73 to A1:R74.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with its new data.
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C7").Value = 'Ñuble'
$ws.Range("D7").Value = 45050
$ws.Range("E7").Value = 16
$ws.Range("F7").Value = 100112001
$ws.Range("G7").Value = 'Berenjena'
$ws.Range("H7").Value = 'Sin especificar'
$ws.Range("I7").Value = 'Primera'
$ws.Range("J7").Value = 60
$ws.Range("K7").Value = 10000
$ws.Range("L7").Value = 11000
$ws.Range("M7").Value = 10500
$ws.Range("N7").Value = '$/caja 60 unidades'
$ws.Range("O7").Value = 'Región de Arica y Parinacota'
$ws.Range("P7").Value = 175
$ws.Range("Q7").Value = 60
$ws.Range("R7").Value = 'Hortaliza'
